$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new player stat rows (26-49)
$ws.Cells.Item(26,1).Value = "Jerry Tillery"
$ws.Cells.Item(26,2).Value = "Group1"
$ws.Cells.Item(26,3).Value = 0.6666666666666666
$ws.Cells.Item(26,4).Value = 32.66666666666666
$ws.Cells.Item(26,5).Value = 17.33333333333333
$ws.Cells.Item(26,6).Value = 15.33333333333333

$ws.Cells.Item(27,1).Value = "Jerry Tillery"
$ws.Cells.Item(27,2).Value = "Group2"
$ws.Cells.Item(27,3).Value = 0.6666666666666666
$ws.Cells.Item(27,4).Value = 23.66666666666667
$ws.Cells.Item(27,5).Value = 10.88888888888889
$ws.Cells.Item(27,6).Value = 12.77777777777778

$ws.Cells.Item(28,1).Value = "Jerry Tillery"
$ws.Cells.Item(28,2).Value = "Difference"
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = -8.999999999999996
$ws.Cells.Item(28,5).Value = -6.444444444444445
$ws.Cells.Item(28,6).Value = -2.555555555555557

$ws.Cells.Item(29,1).Value = "John Jenkins"
$ws.Cells.Item(29,2).Value = "Group1"
$ws.Cells.Item(29,3).Value = 0.6666666666666666
$ws.Cells.Item(29,4).Value = 23.66666666666667
$ws.Cells.Item(29,5).Value = 12.33333333333333
$ws.Cells.Item(29,6).Value = 11.33333333333333

$ws.Cells.Item(30,1).Value = "John Jenkins"
$ws.Cells.Item(30,2).Value = "Group2"
$ws.Cells.Item(30,3).Value = 2
$ws.Cells.Item(30,4).Value = 42.33333333333334
$ws.Cells.Item(30,5).Value = 15
$ws.Cells.Item(30,6).Value = 27.33333333333333

$ws.Cells.Item(31,1).Value = "John Jenkins"
$ws.Cells.Item(31,2).Value = "Difference"
$ws.Cells.Item(31,3).Value = 1.333333333333333
$ws.Cells.Item(31,4).Value = 18.66666666666667
$ws.Cells.Item(31,5).Value = 2.666666666666666
$ws.Cells.Item(31,6).Value = 16

$ws.Cells.Item(32,1).Value = "Linval Joseph"
$ws.Cells.Item(32,2).Value = "Group1"
$ws.Cells.Item(32,3).Value = 0
$ws.Cells.Item(32,4).Value = 54.33333333333334
$ws.Cells.Item(32,5).Value = 29.66666666666667
$ws.Cells.Item(32,6).Value = 24.66666666666667

$ws.Cells.Item(33,1).Value = "Linval Joseph"
$ws.Cells.Item(33,2).Value = "Group2"
$ws.Cells.Item(33,3).Value = 0
$ws.Cells.Item(33,4).Value = 17
$ws.Cells.Item(33,5).Value = 8
$ws.Cells.Item(33,6).Value = 9

$ws.Cells.Item(34,1).Value = "Linval Joseph"
$ws.Cells.Item(34,2).Value = "Difference"
$ws.Cells.Item(34,3).Value = 0
$ws.Cells.Item(34,4).Value = -37.33333333333334
$ws.Cells.Item(34,5).Value = -21.66666666666667
$ws.Cells.Item(34,6).Value = -15.66666666666667

$ws.Cells.Item(35,1).Value = "Ogbonnia Okoronkwo"
$ws.Cells.Item(35,2).Value = "Group1"
$ws.Cells.Item(35,3).Value = 0.3333333333333333
$ws.Cells.Item(35,4).Value = 11.33333333333333
$ws.Cells.Item(35,5).Value = 6.666666666666667
$ws.Cells.Item(35,6).Value = 4.666666666666667

$ws.Cells.Item(36,1).Value = "Ogbonnia Okoronkwo"
$ws.Cells.Item(36,2).Value = "Group2"
$ws.Cells.Item(36,3).Value = 1
$ws.Cells.Item(36,4).Value = 32.66666666666666
$ws.Cells.Item(36,5).Value = 23
$ws.Cells.Item(36,6).Value = 9.666666666666666

$ws.Cells.Item(37,1).Value = "Ogbonnia Okoronkwo"
$ws.Cells.Item(37,2).Value = "Difference"
$ws.Cells.Item(37,3).Value = 0.6666666666666667
$ws.Cells.Item(37,4).Value = 21.33333333333333
$ws.Cells.Item(37,5).Value = 16.33333333333333
$ws.Cells.Item(37,6).Value = 4.999999999999999

$ws.Cells.Item(38,1).Value = "Cameron Jordan"
$ws.Cells.Item(38,2).Value = "Group1"
$ws.Cells.Item(38,3).Value = 4
$ws.Cells.Item(38,4).Value = 54.33333333333334
$ws.Cells.Item(38,5).Value = 36.33333333333334
$ws.Cells.Item(38,6).Value = 18

$ws.Cells.Item(39,1).Value = "Cameron Jordan"
$ws.Cells.Item(39,2).Value = "Group2"
$ws.Cells.Item(39,3).Value = 3
$ws.Cells.Item(39,4).Value = 47.66666666666666
$ws.Cells.Item(39,5).Value = 26
$ws.Cells.Item(39,6).Value = 21.66666666666667

$ws.Cells.Item(40,1).Value = "Cameron Jordan"
$ws.Cells.Item(40,2).Value = "Difference"
$ws.Cells.Item(40,3).Value = -1
$ws.Cells.Item(40,4).Value = -6.666666666666671
$ws.Cells.Item(40,5).Value = -10.33333333333334
$ws.Cells.Item(40,6).Value = 3.666666666666668

$ws.Cells.Item(41,1).Value = "Carl Granderson"
$ws.Cells.Item(41,2).Value = "Group1"
$ws.Cells.Item(41,3).Value = 0.3333333333333333
$ws.Cells.Item(41,4).Value = 16.66666666666667
$ws.Cells.Item(41,5).Value = 10
$ws.Cells.Item(41,6).Value = 6.666666666666667

$ws.Cells.Item(42,1).Value = "Carl Granderson"
$ws.Cells.Item(42,2).Value = "Group2"
$ws.Cells.Item(42,3).Value = 1.666666666666667
$ws.Cells.Item(42,4).Value = 64
$ws.Cells.Item(42,5).Value = 36.66666666666666
$ws.Cells.Item(42,6).Value = 27.33333333333333

$ws.Cells.Item(43,1).Value = "Carl Granderson"
$ws.Cells.Item(43,2).Value = "Difference"
$ws.Cells.Item(43,3).Value = 1.333333333333333
$ws.Cells.Item(43,4).Value = 47.33333333333333
$ws.Cells.Item(43,5).Value = 26.66666666666666
$ws.Cells.Item(43,6).Value = 20.66666666666666

$ws.Cells.Item(44,1).Value = "Deatrich Wise"
$ws.Cells.Item(44,2).Value = "Group1"
$ws.Cells.Item(44,3).Value = 1.333333333333333
$ws.Cells.Item(44,4).Value = 38.66666666666666
$ws.Cells.Item(44,5).Value = 25
$ws.Cells.Item(44,6).Value = 13.66666666666667

$ws.Cells.Item(45,1).Value = "Deatrich Wise"
$ws.Cells.Item(45,2).Value = "Group2"
$ws.Cells.Item(45,3).Value = 2.333333333333333
$ws.Cells.Item(45,4).Value = 47.33333333333334
$ws.Cells.Item(45,5).Value = 23.66666666666667
$ws.Cells.Item(45,6).Value = 23.66666666666667

$ws.Cells.Item(46,1).Value = "Deatrich Wise"
$ws.Cells.Item(46,2).Value = "Difference"
$ws.Cells.Item(46,3).Value = 1
$ws.Cells.Item(46,4).Value = 8.666666666666671
$ws.Cells.Item(46,5).Value = -1.333333333333332
$ws.Cells.Item(46,6).Value = 10

$ws.Cells.Item(47,1).Value = "Dexter Lawrence"
$ws.Cells.Item(47,2).Value = "Group1"
$ws.Cells.Item(47,3).Value = 1.666666666666667
$ws.Cells.Item(47,4).Value = 48.33333333333334
$ws.Cells.Item(47,5).Value = 27
$ws.Cells.Item(47,6).Value = 21.33333333333333

$ws.Cells.Item(48,1).Value = "Dexter Lawrence"
$ws.Cells.Item(48,2).Value = "Group2"
$ws.Cells.Item(48,3).Value = 2
$ws.Cells.Item(48,4).Value = 55
$ws.Cells.Item(48,5).Value = 30
$ws.Cells.Item(48,6).Value = 25

$ws.Cells.Item(49,1).Value = "Dexter Lawrence"
$ws.Cells.Item(49,2).Value = "Difference"
$ws.Cells.Item(49,3).Value = 0.3333333333333333
$ws.Cells.Item(49,4).Value = 6.666666666666664
$ws.Cells.Item(49,5).Value = 3
$ws.Cells.Item(49,6).Value = 3.666666666666668

# Apply alternating row-group fill styles to match existing pattern,
# by copying formats from the existing style-2 / style-3 blocks.
$ws.Range("A2:F4").Copy()
$ws.Range("A26:F28").PasteSpecial(-4122)
$ws.Range("A5:F7").Copy()
$ws.Range("A29:F31").PasteSpecial(-4122)
$ws.Range("A2:F4").Copy()
$ws.Range("A32:F34").PasteSpecial(-4122)
$ws.Range("A5:F7").Copy()
$ws.Range("A35:F37").PasteSpecial(-4122)
$ws.Range("A2:F4").Copy()
$ws.Range("A38:F40").PasteSpecial(-4122)
$ws.Range("A5:F7").Copy()
$ws.Range("A41:F43").PasteSpecial(-4122)
$ws.Range("A2:F4").Copy()
$ws.Range("A44:F46").PasteSpecial(-4122)
$ws.Range("A5:F7").Copy()
$ws.Range("A47:F49").PasteSpecial(-4122)

$excel.CutCopyMode = 0
